# Application 10 Jobs 11/11
# Adds three new job-application rows (Dow Jones, Microsoft, Refonte Learning)
# to the tracking sheet, wires up their "Application Link" hyperlinks, marks
# the previous last link (T-Mobile, C7) as a followed/visited link, wraps the
# long company name in the new row, widens a couple of columns, and leaves
# an extra formatted (but empty) row underneath like the source workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 8 - Dow Jones
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Dow Jones"
$ws.Range("B8").Value = "Data Analyst Intern"
$ws.Range("C8").Value = "https://dowjones.wd1.myworkdayjobs.com/en-US/Dow_Jones_Career/userHome"
$ws.Range("D8").Value = "Data reporting, trend analysis, build dashboards, support data pipeline, visualize data insights"
$ws.Range("E8").Value = "SQL, Python, Excel (advanced formulas), Tableau/Google Analytics, quantitative analysis, cloud (AWS/S3 a plus)"

# ---------------------------------------------------------------------
# Row 9 - Microsoft
# ---------------------------------------------------------------------
$ws.Range("A9").Value = "Microsoft"
$ws.Range("B9").Value = "Explorer Intern - First-Year Students"
$ws.Range("C9").Value = "https://jobs.careers.microsoft.com/global/en/apply/thankyou?Job_id=1782349"
$ws.Range("D9").Value = "Hands-on with development tools, collaborate on design and implementation, experience in Software/Product/Program Mgmt"
$ws.Range("E9").Value = "Pursuing Bachelor’s in tech field, Intro to CS course, calculus, interest in CS/Software Engineering, teamwork, problem-solving"

# ---------------------------------------------------------------------
# Row 10 - Refonte Learning
# ---------------------------------------------------------------------
$ws.Range("A10").Value = "Refonte Learning"
$ws.Range("B10").Value = "AI Engineering Training & Internship"
$ws.Range("C10").Value = "https://www.linkedin.com/jobs/view/4074161641/?alternateChannel=search&refId=W%2BLRz14TtQDBlwjGM%2FUzdQ%3D%3D&trackingId=G8cGFxuowHCancZVCRKAPg%3D%3D"
$ws.Range("D10").Value = "Master AI fundamentals, develop and optimize models, preprocess data, and work on real-world AI projects"
$ws.Range("E10").Value = "TensorFlow, Keras, PyTorch, Scikit-learn; model development, data preprocessing, machine learning (ML, DL, CV); cloud (AWS, GCP)"

# Long company name wraps onto two lines, so the row grows taller.
$ws.Range("A10").WrapText = $true
$ws.Rows.Item(10).RowHeight = 29

# ---------------------------------------------------------------------
# Extra trailing formatted (but empty) row, matching the source sheet.
# ---------------------------------------------------------------------
$ws.Range("A11").WrapText = $true

# ---------------------------------------------------------------------
# Hyperlinks for the "Application Link" column - added in the same order
# they appear in the source workbook's relationship table (Microsoft,
# Dow Jones, then Refonte Learning).
# ---------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("C9"), "https://jobs.careers.microsoft.com/global/en/apply/thankyou?Job_id=1782349", "", "", "https://jobs.careers.microsoft.com/global/en/apply/thankyou?Job_id=1782349")
$ws.Hyperlinks.Add($ws.Range("C8"), "https://dowjones.wd1.myworkdayjobs.com/en-US/Dow_Jones_Career/userHome", "", "", "https://dowjones.wd1.myworkdayjobs.com/en-US/Dow_Jones_Career/userHome")
$ws.Hyperlinks.Add($ws.Range("C10"), "https://www.linkedin.com/jobs/view/4074161641/?alternateChannel=search&refId=W%2BLRz14TtQDBlwjGM%2FUzdQ%3D%3D&trackingId=G8cGFxuowHCancZVCRKAPg%3D%3D", "", "", "https://www.linkedin.com/jobs/view/4074161641/?alternateChannel=search&refId=W%2BLRz14TtQDBlwjGM%2FUzdQ%3D%3D&trackingId=G8cGFxuowHCancZVCRKAPg%3D%3D")

# The previously-last link (T-Mobile, row 7) now reads as a followed /
# visited hyperlink (purple, underlined) now that newer links exist below it.
$ws.Range("C7").Font.Underline = $true
$ws.Range("C7").Font.Color = 8388736

# ---------------------------------------------------------------------
# Column widths tuned to fit the new content.
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 14
$ws.Columns.Item(4).ColumnWidth = 53.66

# ---------------------------------------------------------------------
# Restore the selection to where the editor last left the cursor.
# ---------------------------------------------------------------------
$null = $ws.Range("I16").Select()
